$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (changed) date column C for rows 2-9: 45208 -> 45212
foreach ($r in 2..9) {
    $ws.Cells.Item($r, 3).Value = 45212
}

# Update hyperlink formulas in row 2 (S2, T2, V2, W2, X2, Y2) to include extra filename suffixes
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0484/artfynd/A 31987-2023 artfynd.xlsx", "A 31987-2023")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0484/kartor/A 31987-2023 karta.png", "A 31987-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0484/klagomål/A 31987-2023 fsc-klagomål.docx", "A 31987-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0484/klagomålsmail/A 31987-2023 fsc-klagomål mail.docx", "A 31987-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0484/tillsyn/A 31987-2023 tillsynsbegäran.docx", "A 31987-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0484/ti,llsynsmail/A 31987-2023 tillsynsbegäran mail.docx", "A 31987-2023")'
